$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first data block
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14585
$ws1.Range("F4").Value = 700
$ws1.Range("F6").Value = 603
$ws1.Range("F7").Value = 1548
$ws1.Range("F8").Value = 146

# Sheet "全部类型" (All Types) - second data block
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14585
$ws4.Range("F4").Value = 700
$ws4.Range("F8").Value = 603
$ws4.Range("F9").Value = 1548
$ws4.Range("F11").Value = 146
